# SF Refresh update (except EventExpense):
# On the "Contacts" sheet, a new contact "Adam Taylor" is inserted ahead of
# the existing row-2 data: the former A2 value ("Barry Booth") is shifted
# over to C2, and A2 now holds the new contact name "Adam Taylor". B2
# ("Adam Cole") is left untouched. Finally, the sheet's active selection
# moves from A2:B2 to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

if ($ws.Name -ne "Contacts") {
    $ws = $wb.Worksheets.Item("Contacts")
}

# Shift the existing contact name from A2 into C2 (use Value2 - Value's
# getter isn't reliable for round-tripping text through this bridge).
$ws.Range("C2").Value = $ws.Range("A2").Value2

# Put the new contact in A2.
$ws.Range("A2").Value = "Adam Taylor"

# Match the saved selection state after the edit.
[void]$ws.Range("C8").Select()
